$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D4").Value = "error_occured2"
$ws.Range("D8").Value = "error_occured6"
$ws.Range("D29").Value = "error_occured27"
$ws.Range("D35").Value = "error_occured33"
$ws.Range("D36").Value = "error_occured34"
$ws.Range("D37").Value = "error_occured35"
$ws.Range("D40").Value = "error_occured38"
